$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the displayed recipient e-mail text (leaves the existing hyperlink as-is)
$ws.Range("B1").Value = "bubu@send22u.info"

# Update quantities
$ws.Range("B3").Value = 100
$ws.Range("B4").Value = 500
$ws.Range("B6").Value = 150

# Clear out the now unused barcode column
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D6").ClearContents()

# Reset selection to match the saved view state
$ws.Range("B1").Select()
